# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the
# "e45f33ef-bc4f-4c35-ab14-14de017a65db.md" source file has been
# handed off for translation (zh-cn and de-de) with a new priority
# and updated handoff timestamps.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E3").Value = "Ready for handoff"
$ovw.Range("F3").Value = "Ready for handoff"
$ovw.Range("G3").Value = "2016-09-01 16:17:35"

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-09-01 16:17:30"

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-09-01 16:17:35"

# ---- Column width updates (status columns widened to fit "Ready for handoff") ----
$ovw.Columns.Item(5).ColumnWidth = 16.3
$ovw.Columns.Item(6).ColumnWidth = 16.3
$zhcn.Columns.Item(3).ColumnWidth = 16.3
$dede.Columns.Item(3).ColumnWidth = 16.3
